$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A leading apostrophe is Excel's text-qualifier prefix: it forces a numeric-looking
# string (e.g. "540.62") to be stored as text instead of being parsed into a number,
# matching the inline-string cells used for the Price/Volume columns in this sheet.
$quote = "'"

$ws.Range("D2").Value = "61.769.27"
$ws.Range("E2").Value = "  -4.31%  "
$ws.Range("D3").Value = "2.981.60"
$ws.Range("E3").Value = "  -5.13%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = $quote + "540.62"
$ws.Range("E5").Value = "  -5.86%  "
$ws.Range("D6").Value = $quote + "151.84"
$ws.Range("E6").Value = "  -7.85%  "
$ws.Range("D7").Value = $quote + "0.999"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -1.37%  "
$ws.Range("D9").Value = "2.993.34"
$ws.Range("E9").Value = "  -5.26%  "
$ws.Range("E10").Value = "  -4.07%  "
$ws.Range("D11").Value = $quote + "6.16"
$ws.Range("E11").Value = "  -7.32%  "
$ws.Range("E12").Value = "  -4.18%  "
$ws.Range("D13").Value = "3.501.37"
$ws.Range("E13").Value = "  -5.13%  "
$ws.Range("E14").Value = "  -2.38%  "
$ws.Range("D15").Value = "61.765.04"
$ws.Range("D16").Value = $quote + "23.91"
$ws.Range("E16").Value = "  -4.38%  "
$ws.Range("D17").Value = "2.986.46"
$ws.Range("E17").Value = "  -5.19%  "
$ws.Range("E18").Value = "  -5.82%  "
$ws.Range("D19").Value = $quote + "5.16"
$ws.Range("E19").Value = "  -1.32%  "
$ws.Range("D20").Value = $quote + "12.04"
$ws.Range("E20").Value = "  -4.08%  "
$ws.Range("D21").Value = $quote + "381.07"
$ws.Range("E21").Value = "  -6.74%  "
$ws.Range("E22").Value = "  -5.19%  "
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("D24").Value = $quote + "65.99"
$ws.Range("E24").Value = "  -4.29%  "
$ws.Range("D25").Value = $quote + "0.471"
$ws.Range("E25").Value = "  -2.90%  "
$ws.Range("D26").Value = "3.105.30"
$ws.Range("E26").Value = "  -5.30%  "
$ws.Range("E27").Value = "  -2.74%  "
$ws.Range("E28").Value = "  +0.39%  "
$ws.Range("D29").Value = "0.0₃0941"
$ws.Range("E29").Value = "  -7.77%  "
$ws.Range("D30").Value = $quote + "8.11"
$ws.Range("E30").Value = "  -9.05%  "
$ws.Range("D32").Value = $quote + "20.46"
$ws.Range("E32").Value = "  -3.77%  "
$ws.Range("E33").Value = "  -5.46%  "
$ws.Range("D34").Value = $quote + "159.77"
$ws.Range("E34").Value = "  -2.24%  "
$ws.Range("D35").Value = $quote + "5.92"
$ws.Range("E35").Value = "  -6.04%  "
$ws.Range("E36").Value = "  -6.59%  "
$ws.Range("E37").Value = "  -5.46%  "
$ws.Range("E38").Value = "  -6.70%  "
$ws.Range("E39").Value = "  -8.79%  "
$ws.Range("D40").Value = $quote + "37.57"
$ws.Range("E40").Value = "  -1.69%  "
$ws.Range("D41").Value = "2.421.47"
$ws.Range("E41").Value = "  -8.19%  "
$ws.Range("E42").Value = "  -4.92%  "
$ws.Range("D43").Value = $quote + "21.99"
$ws.Range("E43").Value = "  -7.81%  "
$ws.Range("E44").Value = "  -2.91%  "
$ws.Range("D45").Value = $quote + "0.0589"
$ws.Range("E45").Value = "  -3.91%  "
$ws.Range("D46").Value = $quote + "5.17"
$ws.Range("E46").Value = "  -3.24%  "
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("D48").Value = $quote + "0.0244"
$ws.Range("E48").Value = "  -4.16%  "
$ws.Range("D49").Value = $quote + "19.81"
$ws.Range("E49").Value = "  -7.41%  "
$ws.Range("D50").Value = $quote + "0.0952"
$ws.Range("E50").Value = "  -2.40%  "
$ws.Range("D51").Value = $quote + "266.52"
$ws.Range("E51").Value = "  -8.61%  "
